# Append 4 new rows (n17..n20) to the manifest worksheet, rows 18-21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("n17", "n17_e57_321_1_1.jpeg", "'True", "no_meltpatch", "negative"),
    @("n18", "n18_e58_321_1_2.jpeg", "'True", "no_meltpatch", "negative"),
    @("n19", "n19_e61_321_2_1.jpeg", "'True", "no_meltpatch", "negative"),
    @("n20", "n20_e62_321_2_0.jpeg", "'True", "no_meltpatch", "negative")
)

# Leading apostrophe forces Excel to store the word "True" as text (matching
# the existing rows), rather than auto-converting it to a Boolean value.
$startRow = 18
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $rowData[$c - 1]
    }
}
